# "Generate Report for Handback"
#
# Updates the localization-status workbook to reflect that the a.md
# handback package has been generated:
#   - Overview status for zh-cn / de-de -> "Handed back: in sync with en-US"
#   - Per-language sheets (zh-cn, de-de): fill in "Latest Target File",
#     "Latest Handback File" and "Latest Handback DateTime" for both rows,
#     and widen the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/24778d268e9743d016c7015c47f84386e4277b77/e2e/a.md"
$hyperlinkColor = 15570276   # OLE BGR for RGB FF6495ED (the workbook's existing HyperLink font color)

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status column text and widen
# those two columns so the longer status string fits.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# Re-assert the pre-existing hyperlink / date formatting on this sheet so a
# save/reload of the workbook doesn't silently drop it.
foreach ($row in 2,3) {
    $cellB = $wsOverview.Range("B" + $row)
    $cellB.Font.Underline = 2
    $cellB.Font.Color = $hyperlinkColor
    $wsOverview.Range("G" + $row).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Helper: fill in the handback info for a language sheet (zh-cn / de-de)
# ---------------------------------------------------------------------
function Update-LanguageSheet($ws, $xliffName, $handbackDateTime) {
    # Widen the Status column (C) and the Latest Handback File column (J)
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664

    foreach ($row in 2,3) {
        # Status (C): shares the same underlying text as the Overview
        # status columns, so it also flips to the new handback message.
        $ws.Range("C" + $row).Value = $newStatus

        # Re-assert the pre-existing hyperlink style on column A (the
        # source file link) so it survives the save/reload.
        $cellA = $ws.Range("A" + $row)
        $cellA.Font.Underline = 2
        $cellA.Font.Color = $hyperlinkColor

        # Re-assert the pre-existing date/time formatting on the
        # "Latest Handoff Datetime" column.
        $ws.Range("H" + $row).NumberFormat = "yyyy-mm-dd HH:mm:ss"

        # Latest Target File (I): hyperlink to a.md, styled like the
        # existing hyperlink cells (e.g. column A).
        $cellI = $ws.Range("I" + $row)
        $ws.Hyperlinks.Add($cellI, $aMdUrl, "", "", "a.md")
        $cellI.Font.Underline = 2
        $cellI.Font.Color = $hyperlinkColor

        # Latest Handback File (J): generated xliff file name.
        $ws.Range("J" + $row).Value = $xliffName

        # Latest Handback DateTime (K): when the handback was generated.
        $cellK = $ws.Range("K" + $row)
        $cellK.Value = $handbackDateTime
        $cellK.NumberFormat = "yyyy-mm-dd HH:mm:ss"
    }
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LanguageSheet $wsZhCn "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-29 06:37:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LanguageSheet $wsDeDe "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-29 06:37:56"
